$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the bold/centered/bordered header style from row 1 (A1:Q1) and
# clear out the stray "Unnamed: 0" label that used to sit in A1.
$ws.Range("A1:Q1").ClearFormats()
$ws.Range("A1").ClearContents()

# Corrected data cleaning for pre/post/total fixation data.

# Row 3 - Revisit count
$ws.Range("C3").Value = 31
$ws.Range("D3").Value = 12
$ws.Range("K3").Value = 7
$ws.Range("L3").Value = 28

# Row 4 - Fixation count
$ws.Range("C4").Value = 89
$ws.Range("D4").Value = 20
$ws.Range("K4").Value = 10
$ws.Range("L4").Value = 109

# Row 5 - Dwell time (ms)
$ws.Range("C5").Value = 34794.3
$ws.Range("D5").Value = 9309.56
$ws.Range("K5").Value = 2928.19
$ws.Range("L5").Value = 39032.3

# Row 6 - Dwell time (%)
$ws.Range("B6").Value = 0.49
$ws.Range("C6").Value = 20.27
$ws.Range("D6").Value = 5.42
$ws.Range("E6").Value = 1.88
$ws.Range("H6").Value = 0.49
$ws.Range("I6").Value = 6.04
$ws.Range("J6").Value = 3.38
$ws.Range("K6").Value = 1.71
$ws.Range("L6").Value = 22.74
$ws.Range("M6").Value = 0.39
$ws.Range("N6").Value = 0.55
$ws.Range("O6").Value = 0.13

# Row 7 - Fixation duration (ms)
$ws.Range("C7").Value = 390.95
$ws.Range("D7").Value = 465.48
$ws.Range("K7").Value = 292.82
$ws.Range("L7").Value = 358.09

$wb.Save()
